$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31; this pushes the existing rows 31-36
# down to 32-37, carrying their data (and the D-column date style) along.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with the new weekly price-report entry.
$ws.Range("A31").Value = 4
$ws.Range("B31").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C31").Value = "Los Lagos"
$ws.Range("D31").Value = 44918
$ws.Range("E31").Value = 10
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100101
$ws.Range("H31").Value = "Berries"
$ws.Range("I31").Value = 100101001
$ws.Range("J31").Value = "Arándano (blue)"
$ws.Range("K31").Value = "Sin especificar"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 600
$ws.Range("N31").Value = 3000
$ws.Range("O31").Value = 3500
$ws.Range("P31").Value = 3250
$ws.Range("Q31").Value = "$/bandeja 2 kilos"
$ws.Range("R31").Value = "Provincia de Curicó"
$ws.Range("S31").Value = 1625
$ws.Range("T31").Value = 2
